$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Relocate the "_GoBack" bookmark: it currently sits in the title
#    ("RFA2" | _GoBack | ".1 - ...") and must move to the end of the
#    "4a Il sistema verifica..." paragraph (the location of the last
#    text edit performed below).
# ---------------------------------------------------------------------

# Remove the old bookmark by letting a no-op replace pass over its
# location (the only way this host exposes to drop an existing
# "_GoBack" bookmark).
$old = $d.Content
$old.Find.Execute("RFA2.1", $true, $false, $false, $false, $false, $true, 1, $false, "RFA2.1", 2)

# ---------------------------------------------------------------------
# 2. "Il Cliente inserisce..." sentence: drop the trailing period and
#    add the data-dictionary clause.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Il Cliente inserisce la nuova quantità del prodotto all’interno del carrello.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Il Cliente inserisce la nuova quantità del prodotto all’interno del carrello, in accordo con i criteri del dizionario dei dati, sezione DD_QtOrd",
    2
)

# ---------------------------------------------------------------------
# 3. "4a Il sistema verifica..." sentence: replace "negativo" with the
#    data-dictionary clause.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "4a Il sistema verifica che il cliente ha inserito un valore negativo",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "4a Il sistema verifica che il cliente ha inserito un valore errato, facendo riferimento al dizionario dei dati, sezione DD_QtOrd",
    2
)

# Re-create the "_GoBack" bookmark right after that sentence (last edit
# location), matching where Word leaves it following the final change.
$target = $d.Content
$target.Find.Execute(
    "4a Il sistema verifica che il cliente ha inserito un valore errato, facendo riferimento al dizionario dei dati, sezione DD_QtOrd",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)
$end = $d.Range($target.End, $target.End)
$d.Bookmarks.Add("_GoBack", $end)
